$d = $word.ActiveDocument

# 1. Remove the last (empty) footnote entirely - this deletes both the
#    footnoteReference run in the body and the <w:footnote> element in
#    footnotes.xml (fixes the "empty note" bug mentioned in the commit).
$d.Footnotes.Item(16).Delete()

# 2. Footnote 34 ("gzung ba'i... pe cin.a") - strip the stray trailing "a".
$fn34 = $d.Footnotes.Item(14)
$fn34.Range.Text = " གཟུང་བའི། སྣར་ཐང་། པེ་ཅིན།"

# 3. Footnote 35 - was just a lone shad "།"; fill in the real note text.
$fn35 = $d.Footnotes.Item(15)
$fn35.Range.Text = " གི། ཞེས་པར་མ་གཞན་ནང་མེད།"
